# Applies the edits described by the commit "Update Prequal Table Schema - 21'03.xlsx"
#
# 1. Rename worksheet "RAWMAT_PRE_BATCH_STATUS (2)" -> "RAWMAT_PREQ_SUMMARY"
# 2. Change the active selection on that sheet from A3:A5 (active cell A3)
#    to the single cell B14
# 3. Re-enter the CONCATENATE formula across F23:F28 as one pass so Excel
#    collapses it back into a shared formula group (t="shared"), matching
#    the pattern already used for F3:F67 / F6 elsewhere on the sheet.

$wb = $excel.ActiveWorkbook

# --- 1. Rename the sheet ------------------------------------------------
$ws = $wb.Worksheets.Item("RAWMAT_PRE_BATCH_STATUS (2)")
$ws.Name = "RAWMAT_PREQ_SUMMARY"

# --- 2. Update the selection --------------------------------------------
$ws.Activate()
[void]$ws.Range("B14").Select()

# --- 3. Rebuild F23:F28 as a shared formula ------------------------------
$ws.Range("F23:F28").Formula = '=CONCATENATE(B23," ",C23, D23, IF(ISBLANK(E23), "", CONCATENATE(" ",E23)), ",")'
